$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.183.78"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "2.616.22"
$ws.Range("E3").Value = "  -0.90%  "
$ws.Range("E4").Value = "  +0.03%  "
$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "590.27"
$c.Style = $origStyle
$ws.Range("E5").Value = "  -1.34%  "
$c = $ws.Range("D6")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "165.82"
$c.Style = $origStyle
$ws.Range("E6").Value = "  -0.74%  "
$ws.Range("E7").Value = "  +0.03%  "
$c = $ws.Range("D8")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.531"
$c.Style = $origStyle
$ws.Range("D9").Value = "2.614.99"
$ws.Range("E9").Value = "  -0.93%  "
$c = $ws.Range("D11")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.160"
$c.Style = $origStyle
$ws.Range("E11").Value = "  +1.40%  "
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("E13").Value = "  -0.53%  "
$c = $ws.Range("D14")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "27.32"
$c.Style = $origStyle
$ws.Range("E14").Value = "  -2.19%  "
$ws.Range("D15").Value = "3.089.39"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("E16").Value = "  -2.36%  "
$ws.Range("D17").Value = "67.159.31"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "2.612.78"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("E19").Value = "  -0.66%  "
$c = $ws.Range("D20")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.80"
$c.Style = $origStyle
$ws.Range("E20").Value = "  -0.74%  "
$c = $ws.Range("D21")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "355.06"
$c.Style = $origStyle
$ws.Range("E21").Value = "  -2.10%  "
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("E23").Value = "  -2.78%  "
$c = $ws.Range("D24")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "10.48"
$c.Style = $origStyle
$ws.Range("E24").Value = "  -4.15%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -4.35%  "
$c = $ws.Range("D27")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "69.15"
$c.Style = $origStyle
$ws.Range("E27").Value = "  -2.38%  "
$ws.Range("D28").Value = "2.751.59"
$ws.Range("E28").Value = "  -1.02%  "
$c = $ws.Range("D29")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = $origStyle
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("E30").Value = "  -2.45%  "
$c = $ws.Range("D31")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "542.89"
$c.Style = $origStyle
$ws.Range("E31").Value = "  -2.01%  "
$c = $ws.Range("D32")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "7.87"
$c.Style = $origStyle
$ws.Range("E32").Value = "  -2.12%  "
$ws.Range("E34").Value = "  -2.74%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("E37").Value = "  -3.52%  "
$c = $ws.Range("D38")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "159.11"
$c.Style = $origStyle
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("E40").Value = "  -2.18%  "
$c = $ws.Range("D41")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "18.25"
$c.Style = $origStyle
$ws.Range("E41").Value = "  +1.78%  "
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("E43").Value = "  -2.34%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  -4.40%  "
$ws.Range("E46").Value = "  -1.00%  "
$c = $ws.Range("D47")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "151.48"
$c.Style = $origStyle
$ws.Range("E47").Value = "  -1.27%  "
$ws.Range("E48").Value = "  -3.28%  "
$ws.Range("E49").Value = "  -3.02%  "
$c = $ws.Range("D50")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.71"
$c.Style = $origStyle
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("E51").Value = "  -1.27%  "
